$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of new-string introduction matches the target shared-strings table:
# 14 POST, 15 candidate/job/apply, 16 job_id, 17 candidate/job/apply/list,
# 18 candidate/banners, 19 candidate/news, 20 candidate/news-categories,
# 21 candidate/news/{id}

# Row 11: POST | candidate/job/apply | job_id
$ws.Range("A11").Value = "POST"
$ws.Range("C11").Value = "candidate/job/apply"
$ws.Range("D11").Value = "job_id"

# Row 12: GET | candidate/job/apply/list
$ws.Range("A12").Value = "GET"
$ws.Range("C12").Value = "candidate/job/apply/list"

# Row 9: candidate/banners (styled like C8, Consolas red font) - reuse existing style
$ws.Range("C9").Value = "candidate/banners"
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

# Row 10: empty cell but keep same style as C8/C9
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 14: GET | candidate/news
$ws.Range("A14").Value = "GET"
$ws.Range("C14").Value = "candidate/news"

# Row 15: candidate/news-categories
$ws.Range("C15").Value = "candidate/news-categories"

# Row 16: candidate/news/{id}
$ws.Range("C16").Value = "candidate/news/{id}"

# Row 13 intentionally left blank

# Update selection to match final saved state
$ws.Range("A14").Select()
